# "updated api for bdd"
# - Rename the existing "Sheet1" to "RestAssured".
# - Add a new worksheet "CreateUser" after it, with a small status/username/
#   password/email table and a "v@v.com" mailto hyperlink.
# - Leave RestAssured's own selection on C28, and finish with CreateUser as
#   the active (selected) tab, selection on B2.

$wb = $excel.ActiveWorkbook

# --- Rename the first sheet --------------------------------------------
$restAssured = $wb.Worksheets.Item(1)
$restAssured.Name = "RestAssured"

# Park the selection on RestAssured where it was left before switching
# focus to the new sheet (matches the saved view state).
[void]$restAssured.Range("C28").Select()

# --- Add the new sheet after RestAssured --------------------------------
$createUser = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$createUser.Name = "CreateUser"

# Header row. Populate columns B:E first so the shared-string table picks
# up "status"/"username"/"password"/"email" (etc.) ahead of "User".
$createUser.Range("B1").Value = "status"
$createUser.Range("C1").Value = "username"
$createUser.Range("D1").Value = "password"
$createUser.Range("E1").Value = "email"

# Data row.
$createUser.Range("B2").Value = "OK"
$createUser.Range("C2").Value = "vipin"
$createUser.Range("D2").Value = "password"
$createUser.Range("E2").Value = "v@v.com"

# Column A labels both rows as "User" (added last so it lands at the end
# of the shared-string table).
$createUser.Range("A1").Value = "User"
$createUser.Range("A2").Value = "User"

# Turn the e-mail cell into a live mailto: link (applies the built-in
# Hyperlink style/font automatically).
$createUser.Hyperlinks.Add($createUser.Range("E2"), "mailto:v@v.com")

# Fit column A to its content.
$createUser.Columns.Item(1).AutoFit()

# Final selection/active tab: CreateUser, cell B2.
[void]$createUser.Range("B2").Select()
